$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin rows/values per upstream data refresh.
# Columns B (Coin), C (Link), D (Price), E (Volume(1h)) are plain text cells;
# force text number format so values like percentages/decimals are not
# reinterpreted as numeric/percentage types by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '308.48'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.38%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '40.71'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-0.66%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.125'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.53%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07620'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.01%'
$ws.Range('B6').NumberFormat = '@'
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.608'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.73%'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'BTSEToken'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '2.456'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.53%'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9077'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '0.13%'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1275'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '28.93%'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1798'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '1.76%'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.09024'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-1.19%'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.04326'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-1.88%'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.1043'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.92%'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001254'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.49%'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.005843'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.58%'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.346'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.61%'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.279'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.06%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.991'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.33%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.1393'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '3.34%'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-4.75%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.04045'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-2.89%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.001274'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '4.97%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.004039'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.72%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0001273'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-2.08%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '24.81%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02423'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '0.23%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05214'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '1.44%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.007853'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.29%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1298'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.79%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.006802'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-3.75%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.001900'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-2.48%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007389'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-11.72%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3348'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '9.79%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006899'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '8.20%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000751'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.16%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1076'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '1,348.79%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002104'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.16%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002003'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.16%'
